# Apply the Apr 13 2023 cryptos refresh to the active worksheet.
# Column layout: A=index(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
# All D/E cells hold text (not numeric) values in the source data, so we
# temporarily force text format while writing to stop Excel's automatic
# "looks like a number" conversion, then clear the formatting again so the
# cells end up with no explicit style, just like the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Rows 37 and 38 swapped their coin (Hedera <-> InternetComputer(DFINITY))
# along with all associated data.
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "5.478"
$ws.Range("E37").Value = "  +1.00%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06671"
$ws.Range("E38").Value = "  +3.85%  "

# Refreshed price / volume figures for all other rows.
$ws.Range("D2").Value = "30.397.69"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.013.25"
$ws.Range("E3").Value = "  +4.77%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "324.97"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.5136"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").Value = "0.4261"
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("D9").Value = "0.08768"
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("D10").Value = "43.52"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").Value = "2.019.68"
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").Value = "6.659"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "7.471"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "94.24"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "0.06542"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "18.88"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "6.210"
$ws.Range("E22").Value = "  +4.37%  "
$ws.Range("D23").Value = "30.463.06"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("D25").Value = "2.251"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "2.252.47"
$ws.Range("E26").Value = "  +4.93%  "
$ws.Range("D27").Value = "22.46"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "162.10"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "2.451"
$ws.Range("E29").Value = "  +4.82%  "
$ws.Range("D30").Value = "131.39"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "0.1055"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").Value = "6.088"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "3.830"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "1.368"
$ws.Range("E35").Value = "  +14.43%  "
$ws.Range("D36").Value = "0.02533"
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D39").Value = "12.42"
$ws.Range("E39").Value = "  +8.95%  "
$ws.Range("D40").Value = "9.223"
$ws.Range("E40").Value = "  +5.26%  "
$ws.Range("D41").Value = "0.2218"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "1.241"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "13.67"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").Value = "0.6181"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "3.639"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "1.259"
$ws.Range("E49").Value = "  +4.04%  "
$ws.Range("D50").Value = "125.21"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "81.29"
$ws.Range("E51").Value = "  +2.80%  "

# Remove the temporary text formatting so the cells end up with no explicit
# style, matching the original workbook.
$dataRange.ClearFormats()
